# Generate Report for Handoff
# Adds two new files (324a4567-38de-4d2d-975e-9d6a52a3674c and
# 5c4c6826-7756-4723-a923-e65d0f2de573) to the Overview / zh-cn / de-de
# report sheets, each appearing as "Ready for handoff".

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b8dca5e31295c1e8f1d2eff6a350d61283c92c6f/e2e/"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.md"
$wsOverview.Range("B4").Value = "e2e\324a4567-38de-4d2d-975e-9d6a52a3674c.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Ready for handoff"
$wsOverview.Range("F4").Value = "Ready for handoff"
$wsOverview.Range("G4").Value = "2016-10-19 16:38:34"

$wsOverview.Range("A5").Value = "5c4c6826-7756-4723-a923-e65d0f2de573.md"
$wsOverview.Range("B5").Value = "e2e\5c4c6826-7756-4723-a923-e65d0f2de573.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-10-19 16:38:34"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), ($repoBase + "324a4567-38de-4d2d-975e-9d6a52a3674c.md"), [Type]::Missing, [Type]::Missing, "e2e\324a4567-38de-4d2d-975e-9d6a52a3674c.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), ($repoBase + "5c4c6826-7756-4723-a923-e65d0f2de573.md"), [Type]::Missing, [Type]::Missing, "e2e\5c4c6826-7756-4723-a923-e65d0f2de573.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.md"
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Ready for handoff"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "False"
$wsZh.Range("G4").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.5e02424a11a8004174b34e7fcb9bc4a1b236430b.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-10-19 16:38:22"
$wsZh.Range("K4").Value = "0001-01-01 00:00:00"
$wsZh.Range("M4").Value = "True"
$wsZh.Range("O4").Value = "False"

$wsZh.Range("A5").Value = "5c4c6826-7756-4723-a923-e65d0f2de573.md"
$wsZh.Range("B5").Value = ".md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("D5").Value = "e2e"
$wsZh.Range("E5").Value = "ht"
$wsZh.Range("F5").Value = "False"
$wsZh.Range("G5").Value = "5c4c6826-7756-4723-a923-e65d0f2de573.ebd51b133f7dc726b1cf4177bf4407a103b4a5db.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-10-19 16:38:22"
$wsZh.Range("K5").Value = "0001-01-01 00:00:00"
$wsZh.Range("M5").Value = "True"
$wsZh.Range("O5").Value = "False"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), ($repoBase + "324a4567-38de-4d2d-975e-9d6a52a3674c.md"), [Type]::Missing, [Type]::Missing, "324a4567-38de-4d2d-975e-9d6a52a3674c.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), ($repoBase + "5c4c6826-7756-4723-a923-e65d0f2de573.md"), [Type]::Missing, [Type]::Missing, "5c4c6826-7756-4723-a923-e65d0f2de573.md")

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.md"
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Ready for handoff"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "False"
$wsDe.Range("G4").Value = "324a4567-38de-4d2d-975e-9d6a52a3674c.5e02424a11a8004174b34e7fcb9bc4a1b236430b.de-de.xlf"
$wsDe.Range("H4").Value = "2016-10-19 16:38:34"
$wsDe.Range("K4").Value = "0001-01-01 00:00:00"
$wsDe.Range("M4").Value = "True"
$wsDe.Range("O4").Value = "False"

$wsDe.Range("A5").Value = "5c4c6826-7756-4723-a923-e65d0f2de573.md"
$wsDe.Range("B5").Value = ".md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("D5").Value = "e2e"
$wsDe.Range("E5").Value = "ht"
$wsDe.Range("F5").Value = "False"
$wsDe.Range("G5").Value = "5c4c6826-7756-4723-a923-e65d0f2de573.ebd51b133f7dc726b1cf4177bf4407a103b4a5db.de-de.xlf"
$wsDe.Range("H5").Value = "2016-10-19 16:38:34"
$wsDe.Range("K5").Value = "0001-01-01 00:00:00"
$wsDe.Range("M5").Value = "True"
$wsDe.Range("O5").Value = "False"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), ($repoBase + "324a4567-38de-4d2d-975e-9d6a52a3674c.md"), [Type]::Missing, [Type]::Missing, "324a4567-38de-4d2d-975e-9d6a52a3674c.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), ($repoBase + "5c4c6826-7756-4723-a923-e65d0f2de573.md"), [Type]::Missing, [Type]::Missing, "5c4c6826-7756-4723-a923-e65d0f2de573.md")

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Column width adjustments (Status columns widened to fit "Ready for handoff")
# ---------------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZh.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDe.Columns.Item(3).ColumnWidth = 16.333333333333332
